$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.046.17"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "2.344.88"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.57"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.09"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.06"
$ws.Range("E11").Value = "  -6.29%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "2.707.83"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.41"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.78"
$ws.Range("E15").Value = "  -6.28%  "
$ws.Range("D16").Value = "2.374.32"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.759"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").Value = "40.035.32"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "0.0₃0901"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.17"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.62"
$ws.Range("E22").Value = "  -5.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.99"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("E24").Value = "  -5.65%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.64"
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.97"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.76"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.10"
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -6.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0990"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.64"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").Value = "1.970.83"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0264"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.48"
$ws.Range("E45").Value = "  -6.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.52"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -7.41%  "
$ws.Range("D48").Value = "2.569.18"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.10"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.44"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.24"
$ws.Range("E51").Value = "  -3.97%  "
